# Update "paises.xlsx" with latest COVID-19 country & Spain provincias data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 22:05"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1660622
$ws.Range("C4").Value = 15528
$ws.Range("D4").Value = 439487
$ws.Range("E4").Value = 1122704
$ws.Range("G4").Value = 784
$ws.Range("H4").Value = 98431

# Row 5: Brasil
$ws.Range("B5").Value = 341048
$ws.Range("C5").Value = 10158
$ws.Range("E5").Value = 183936
$ws.Range("G5").Value = 634
$ws.Range("H5").Value = 21682

# Row 11: Alemania
$ws.Range("B11").Value = 179986
$ws.Range("C11").Value = 273
$ws.Range("E11").Value = 11720
$ws.Range("G11").Value = 14
$ws.Range("H11").Value = 8366

# Row 54: Barein
$ws.Range("B54").Value = 8802
$ws.Range("C54").Value = 388
$ws.Range("D54").Value = 4465
$ws.Range("E54").Value = 4324

# Row 84: Costa de Marfil
$ws.Range("B84").Value = 2366
$ws.Range("C84").Value = 25
$ws.Range("D84").Value = 1188
$ws.Range("E84").Value = 1148
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 30

# Row 101: Maldivas
$ws.Range("D101").Value = 115
$ws.Range("E101").Value = 1194
